# Apply cryptos-list refresh (prices/volumes updated; ImmutableX/NEARProtocol rows swapped)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 29 & 30: NEARProtocol and ImmutableX swapped places in the ranking ---
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "63.021.25"
$ws.Range("D3").Value = "3.031.05"
$ws.Range("D5").Value = "'594.09"
$ws.Range("D6").Value = "'154.07"
$ws.Range("D8").Value = "3.026.93"
$ws.Range("D14").Value = "'35.63"
$ws.Range("D16").Value = "3.533.63"
$ws.Range("D18").Value = "62.961.50"
$ws.Range("D19").Value = "3.031.27"
$ws.Range("D20").Value = "'453.83"
$ws.Range("D22").Value = "'0.699"
$ws.Range("D24").Value = "'11.53"
$ws.Range("D27").Value = "'12.44"
$ws.Range("D29").Value = "'7.51"
$ws.Range("D30").Value = "'2.28"
$ws.Range("D35").Value = "0.0₃0866"
$ws.Range("D38").Value = "'3.16"
$ws.Range("D41").Value = "'50.46"
$ws.Range("D42").Value = "'9.10"
$ws.Range("D44").Value = "'44.66"
$ws.Range("D45").Value = "'395.04"
$ws.Range("D46").Value = "'0.0361"
$ws.Range("D47").Value = "2.721.51"
$ws.Range("D48").Value = "'132.83"
$ws.Range("D49").Value = "'25.73"

# --- Column E (Volume 1h %) updates ---
$ws.Range("E2").Value = "  +3.20%  "
$ws.Range("E3").Value = "  +2.10%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("E6").Value = "  +9.04%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  +2.22%  "
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("E10").Value = "  +17.36%  "
$ws.Range("E11").Value = "  +4.98%  "
$ws.Range("E12").Value = "  +3.08%  "
$ws.Range("E13").Value = "  +4.26%  "
$ws.Range("E14").Value = "  +5.12%  "
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("E16").Value = "  +2.00%  "
$ws.Range("E17").Value = "  +4.04%  "
$ws.Range("E18").Value = "  +2.84%  "
$ws.Range("E19").Value = "  +1.75%  "
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("E22").Value = "  +3.39%  "
$ws.Range("E23").Value = "  +4.50%  "
$ws.Range("E24").Value = "  +12.72%  "
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("E26").Value = "  +9.70%  "
$ws.Range("E27").Value = "  +4.76%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  +7.56%  "
$ws.Range("E30").Value = "  +12.87%  "
$ws.Range("E31").Value = "  +1.96%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  +2.29%  "
$ws.Range("E34").Value = "  +3.91%  "
$ws.Range("E35").Value = "  +8.59%  "
$ws.Range("E36").Value = "  +3.79%  "
$ws.Range("E37").Value = "  +3.52%  "
$ws.Range("E38").Value = "  +12.57%  "
$ws.Range("E39").Value = "  +8.97%  "
$ws.Range("E40").Value = "  +3.95%  "
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("E43").Value = "  +17.73%  "
$ws.Range("E44").Value = "  +16.82%  "
$ws.Range("E45").Value = "  +2.30%  "
$ws.Range("E46").Value = "  +4.41%  "
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("E48").Value = "  +2.97%  "
$ws.Range("E49").Value = "  +12.09%  "
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("E51").Value = "  +8.50%  "

